$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet has extra gamma-index / QA columns (S:Y) that are being trimmed down
# to just the "MedianDoseDev" column (originally column V). Move column V's
# data+formatting into column S (overwriting the old "Global Max Gamma Index"
# data there), then clear out the now-redundant columns T:Y.
$ws.Range("V1:V58").Copy($ws.Range("S1:S58")) | Out-Null
$ws.Range("T1:Y58").Clear() | Out-Null

# Reflect the selection the author left behind after trimming the columns.
$ws.Activate()
$ws.Range("T1:V1048576").Select() | Out-Null
